# Post session 2 rebalance
# Raw material "Qty" scores are rebalanced down to a flat 2, and the
# crafted-item "Qty" column (previously computed from the raw material
# quantities) is hard-coded to a flat 1 for every craftable resource.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generics")
$ws.Activate()

# Raw materials (rows 3-9) - rebalance Qty column (B)
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2

# Crafted items (rows 11-18) - replace the Qty formulas with a flat value
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1

# Leave the cursor where the author left off
[void]$ws.Range("B19").Select()
